# Applies the "Add files via upload" diff to the israel_ligat-ha-al_2023-2024
# sheet:
#   1) Swap the match-detail columns (F:V) between row 16 <-> row 17
#   2) Swap the match-detail columns (F:V) between row 18 <-> row 19
#   3) Append a new row 28 (Maccabi Haifa vs Sakhnin)
#
# Note: Range.Value getter in this runtime doesn't reliably surface the
# underlying scalar (it stubs out to a placeholder string) - .Formula
# reads back the real text/number for both string and numeric cells, so
# it is used here for every read. Writes use .Value, which works fine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($col, $r1, $r2) {
    $addr1 = "$col$r1"
    $addr2 = "$col$r2"
    $v1 = $ws.Range($addr1).Formula
    $v2 = $ws.Range($addr2).Formula
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}

# Columns F..V, excluding K, O, S (those keep their original per-row value).
$swapCols = @("F","G","H","I","J","L","M","N","P","Q","R","T","U","V")

foreach ($col in $swapCols) {
    Swap-Cell $col 16 17
}

foreach ($col in $swapCols) {
    Swap-Cell $col 18 19
}

# --- Append new row 28: Maccabi Haifa 1-1 Sakhnin -----------------------

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "israel"
$ws.Range("C28").Value = "ligat-ha-al"
$ws.Range("D28").Value = "2023-2024"
$ws.Range("E28").Value = 45196.79166666666
$ws.Range("F28").Value = "Maccabi Haifa"
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = "Sakhnin"
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = 1.17
$ws.Range("K28").Value = "20/09/2023 18:12"
$ws.Range("L28").Value = 1.27
$ws.Range("M28").Value = "27/09/2023 18:54"
$ws.Range("N28").Value = 6.94
$ws.Range("O28").Value = "20/09/2023 18:12"
$ws.Range("P28").Value = 6.23
$ws.Range("Q28").Value = "27/09/2023 18:54"
$ws.Range("R28").Value = 12.12
$ws.Range("S28").Value = "20/09/2023 18:12"
$ws.Range("T28").Value = 9.31
$ws.Range("U28").Value = "27/09/2023 18:54"
$ws.Range("V28").Value = "https://www.betexplorer.com/football/israel/ligat-ha-al/maccabi-haifa-sakhnin/EmhlJEs4/"

# Match the existing column styles used throughout the sheet: column A is
# bold/bordered (s=1), column E is the datetime number format (s=2).
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("E27").Copy()
$ws.Range("E28").PasteSpecial(-4122)
